$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.557.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.081.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +5.39%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.00"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +7.27%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.435"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.29"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.69%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.608.15"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.554.44"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.081.11"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.10"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.04"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.19"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "339.13"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +7.61%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.500"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.62"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.49%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.45%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0947"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +12.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.49"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.46%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.83"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.89%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.42"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.56"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.91"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.08"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.24"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0679"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.120.02"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.98"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.86"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.40%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.672"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.43%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.270.32"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.40"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0251"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.962"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.33"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +8.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.87"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.03%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.42%  "
